$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Disambiguate the duplicate "Centraal Station" stop name (row 19) by
# renaming it to "Centraal Station_B".
$ws.Range("A19").Value = "Centraal Station_B"
